$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9064210071256866
$ws.Range("C2").Value = 0.2203067207800586
$ws.Range("D2").Value = 0.2023424355200092
$ws.Range("E2").Value = 0.1547909353960861
$ws.Range("F2").Value = 1.177518146770083
$ws.Range("I2").Value = 0.5231711738452685
$ws.Range("J2").Value = 0.1596138729011791
$ws.Range("M2").Value = 0.3703977477924099
$ws.Range("O2").Value = 2.689013073750971
$ws.Range("B3").Value = 0.8031573600234765
$ws.Range("C3").Value = 0.1929127848275414
$ws.Range("D3").Value = 0.2000383646797488
$ws.Range("E3").Value = 0.1548333156998396
$ws.Range("F3").Value = 1.183726354629769
$ws.Range("I3").Value = 0.5328619925021876
$ws.Range("J3").Value = 0.1609187147279663
$ws.Range("M3").Value = 0.3430046127083344
$ws.Range("O3").Value = 2.711440572946941
$ws.Range("B4").Value = 0.7396247032718861
$ws.Range("C4").Value = 0.1760380690235195
$ws.Range("D4").Value = 0.1986851779472474
$ws.Range("E4").Value = 0.1549289865395949
$ws.Range("F4").Value = 1.188358400747731
$ws.Range("I4").Value = 0.5392324207284123
$ws.Range("J4").Value = 0.1618099164968676
$ws.Range("M4").Value = 0.326233312942243
$ws.Range("O4").Value = 2.727302318494168
$ws.Range("B5").Value = 0.7137039781729584
$ws.Range("C5").Value = 0.1691481309999858
$ws.Range("D5").Value = 0.1981492872878903
$ws.Range("E5").Value = 0.154985507000454
$ws.Range("F5").Value = 1.190452051768879
$ws.Range("I5").Value = 0.5419339037214872
$ws.Range("J5").Value = 0.1621957163346224
$ws.Range("M5").Value = 0.3194114119666978
$ws.Range("O5").Value = 2.734291136818598
$ws.Range("B6").Value = 0.7093980583691462
$ws.Range("C6").Value = 0.1680032675426162
$ws.Range("D6").Value = 0.1980612440413552
$ws.Range("E6").Value = 0.1549959517686013
$ws.Range("F6").Value = 1.190812140644795
$ws.Range("I6").Value = 0.5423888482487342
$ws.Range("J6").Value = 0.16226114450129
$ws.Range("M6").Value = 0.3182794103298221
$ws.Range("O6").Value = 2.735483311316642
$ws.Range("B7").Value = 0.7392752492501131
$ws.Range("C7").Value = 0.1759452024387542
$ws.Range("D7").Value = 0.1986778877022459
$ws.Range("E7").Value = 0.1549296777738753
$ws.Range("F7").Value = 1.188385802400127
$ws.Range("I7").Value = 0.5392684270313275
$ws.Range("J7").Value = 0.1618150279209978
$ws.Range("M7").Value = 0.326141259066361
$ws.Range("O7").Value = 2.72739444740931
$ws.Range("B8").Value = 0.8708432822335794
$ws.Range("C8").Value = 0.2108728857943731
$ws.Range("D8").Value = 0.2015352660582295
$ws.Range("E8").Value = 0.1547911040532739
$ws.Range("F8").Value = 1.179488394169475
$ws.Range("I8").Value = 0.526425238430333
$ws.Range("J8").Value = 0.1600450923155208
$ws.Range("M8").Value = 0.3609428344195109
$ws.Range("O8").Value = 2.696311561930059
$ws.Range("B9").Value = 1.127770804551631
$ws.Range("C9").Value = 0.2789171793599792
$ws.Range("D9").Value = 0.2076239239023749
$ws.Range("E9").Value = 0.1550712818679507
$ws.Range("F9").Value = 1.168557834434189
$ws.Range("I9").Value = 0.50458221266123
$ws.Range("J9").Value = 0.1572890219107883
$ws.Range("M9").Value = 0.4295561425677903
$ws.Range("O9").Value = 2.65198883944575
$ws.Range("B10").Value = 1.315819253882864
$ws.Range("C10").Value = 0.3286211394539293
$ws.Range("D10").Value = 0.2123900162211214
$ws.Range("E10").Value = 0.1556128750754446
$ws.Range("F10").Value = 1.164514829431809
$ws.Range("I10").Value = 0.4905825250120373
$ws.Range("J10").Value = 0.1557006066864268
$ws.Range("M10").Value = 0.4801756734451672
$ws.Range("O10").Value = 2.629618584400873
$ws.Range("B11").Value = 1.401200570918547
$ws.Range("C11").Value = 0.3511674419926862
$ws.Range("D11").Value = 0.214621186423031
$ws.Range("E11").Value = 0.1559320291363591
$ws.Range("F11").Value = 1.163544680145549
$ws.Range("I11").Value = 0.4846607805398371
$ws.Range("J11").Value = 0.1550729348998345
$ws.Range("M11").Value = 0.5032465122211534
$ws.Range("O11").Value = 2.621667208066299
$ws.Range("B12").Value = 1.433507491900514
$ws.Range("C12").Value = 0.359695562624438
$ws.Range("D12").Value = 0.2154750747442336
$ws.Range("E12").Value = 0.1560633346873175
$ws.Range("F12").Value = 1.163302512095299
$ws.Range("I12").Value = 0.4824828185386529
$ws.Range("J12").Value = 0.1548489124216559
$ws.Range("M12").Value = 0.51198878835865
$ws.Range("O12").Value = 2.618977114009141
$ws.Range("B13").Value = 1.426550763198634
$ws.Range("C13").Value = 0.3578593163812229
$ws.Range("D13").Value = 0.2152907758438118
$ws.Range("E13").Value = 0.1560345912419407
$ws.Range("F13").Value = 1.163349094915262
$ws.Range("I13").Value = 0.4829490105075216
$ws.Range("J13").Value = 0.1548965517444074
$ws.Range("M13").Value = 0.5101057295947697
$ws.Range("O13").Value = 2.619542185470465
$ws.Range("B14").Value = 1.403858994154746
$ws.Range("C14").Value = 0.351869251966491
$ws.Range("D14").Value = 0.2146912565268622
$ws.Range("E14").Value = 0.1559426224094835
$ws.Range("F14").Value = 1.163522246139109
$ws.Range("I14").Value = 0.4844803048381969
$ws.Range("J14").Value = 0.1550542305142528
$ws.Range("M14").Value = 0.5039656298167046
$ws.Range("O14").Value = 2.621439455147652
$ws.Range("B15").Value = 1.389956307120826
$ws.Range("C15").Value = 0.3481988918565548
$ws.Range("D15").Value = 0.2143252025432361
$ws.Range("E15").Value = 0.1558876490892693
$ws.Range("F15").Value = 1.163644618946648
$ws.Range("I15").Value = 0.4854266699988798
$ws.Range("J15").Value = 0.1551525931352735
$ws.Range("M15").Value = 0.5002053902025665
$ws.Range("O15").Value = 2.622643409599306
$ws.Range("B16").Value = 1.310235960142506
$ws.Range("C16").Value = 0.3271463551296279
$ws.Range("D16").Value = 0.2122454664075093
$ws.Range("E16").Value = 0.1555934806909249
$ws.Range("F16").Value = 1.16459573873756
$ws.Range("I16").Value = 0.4909785347968203
$ws.Range("J16").Value = 0.1557435376216567
$ws.Range("M16").Value = 0.4786687825430889
$ws.Range("O16").Value = 2.630183080738959
$ws.Range("B17").Value = 1.261287173421636
$ws.Range("C17").Value = 0.314214534186533
$ws.Range("D17").Value = 0.2109857116947325
$ws.Range("E17").Value = 0.1554316467073846
$ws.Range("F17").Value = 1.16540197188182
$ws.Range("I17").Value = 0.4944990177939559
$ws.Range("J17").Value = 0.1561303819388478
$ws.Range("M17").Value = 0.4654676635006325
$ws.Range("O17").Value = 2.63537898482133
$ws.Range("B18").Value = 1.23311789754564
$ws.Range("C18").Value = 0.306770468213756
$ws.Range("D18").Value = 0.2102670731584624
$ws.Range("E18").Value = 0.1553454155445024
$ws.Range("F18").Value = 1.165947479166562
$ws.Range("I18").Value = 0.496565946332737
$ws.Range("J18").Value = 0.1563618166417378
$ws.Range("M18").Value = 0.4578788741521649
$ws.Range("O18").Value = 2.638576897762817
$ws.Range("B19").Value = 1.223577701637282
$ws.Range("C19").Value = 0.304249013864279
$ws.Range("D19").Value = 0.2100247767982921
$ws.Range("E19").Value = 0.1553173964162937
$ws.Range("F19").Value = 1.166146217576355
$ws.Range("I19").Value = 0.4972729862771175
$ws.Range("J19").Value = 0.1564417100678916
$ws.Range("M19").Value = 0.4553101659446952
$ws.Range("O19").Value = 2.639695585323125
$ws.Range("B20").Value = 1.266499443855139
$ws.Range("C20").Value = 0.3155917752665118
$ws.Range("D20").Value = 0.2111192006019138
$ws.Range("E20").Value = 0.1554481652332846
$ws.Range("F20").Value = 1.16530768097428
$ws.Range("I20").Value = 0.4941199037559514
$ws.Range("J20").Value = 0.1560882771492089
$ws.Range("M20").Value = 0.4668725193656798
$ws.Range("O20").Value = 2.634804196268362
$ws.Range("B21").Value = 1.410524812380743
$ws.Range("C21").Value = 0.353628945080402
$ws.Range("D21").Value = 0.214867106415511
$ws.Range("E21").Value = 0.1559693524412857
$ws.Range("F21").Value = 1.163467987452037
$ws.Range("I21").Value = 0.4840287749566734
$ws.Range("J21").Value = 0.1550075454961153
$ws.Range("M21").Value = 0.5057689706973605
$ws.Range("O21").Value = 2.620873463658285
$ws.Range("B22").Value = 1.504506451993564
$ws.Range("C22").Value = 0.3784318444796781
$ws.Range("D22").Value = 0.2173689610883542
$ws.Range("E22").Value = 0.1563708743886245
$ws.Range("F22").Value = 1.162995550838502
$ws.Range("I22").Value = 0.4778095485603338
$ws.Range("J22").Value = 0.1543808718294173
$ws.Range("M22").Value = 0.5312239465486641
$ws.Range("O22").Value = 2.613639937952513
$ws.Range("B23").Value = 1.454360718269356
$ws.Range("C23").Value = 0.3651993899204058
$ws.Range("D23").Value = 0.2160289056247677
$ws.Range("E23").Value = 0.1561510081649224
$ws.Range("F23").Value = 1.163180833325896
$ws.Range("I23").Value = 0.4810943960192304
$ws.Range("J23").Value = 0.1547080464401489
$ws.Range("M23").Value = 0.5176351907092709
$ws.Range("O23").Value = 2.617329088694362
$ws.Range("B24").Value = 1.264143063377162
$ws.Range("C24").Value = 0.3149691537925321
$ws.Range("D24").Value = 0.2110588327841185
$ws.Range("E24").Value = 0.155440675996207
$ws.Range("F24").Value = 1.165350054506504
$ws.Range("I24").Value = 0.4942911674227979
$ws.Range("J24").Value = 0.1561072845895737
$ws.Range("M24").Value = 0.4662373817858594
$ws.Range("O24").Value = 2.635063401864386
$ws.Range("B25").Value = 1.058386627011942
$ws.Range("C25").Value = 0.2605589639083803
$ws.Range("D25").Value = 0.205925145336991
$ws.Range("E25").Value = 0.1549364739725405
$ws.Range("F25").Value = 1.170815612382199
$ws.Range("I25").Value = 0.5101323877939947
$ws.Range("J25").Value = 0.1579580159245637
$ws.Range("M25").Value = 0.410956613367901
$ws.Range("O25").Value = 2.662193166503272
